$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new row above row 5, pushing the existing rows 5-6 down to 6-7.
$ws.Rows.Item(5).Insert()

# The inserted row inherited formatting/content placeholders from the header
# row (row 4, bold/filled). Clear the whole row so it starts out empty again,
# matching the plain (unstyled) look of the other data rows.
$ws.Rows.Item(5).Clear()

# Re-apply the date number format to column A, matching the other data rows.
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Populate the new transaction (row 5) with the latest contract-note data.
$ws.Range("A5").Value = 46066
$ws.Range("B5").Value = "BSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 866.95
$ws.Range("F5").Value = 4365.27
$ws.Range("G5").Value = "CN#252611910666"
$ws.Range("H5").Value = 4.3787
$ws.Range("I5").Value = 26.1397
$ws.Range("J5").Formula = "=Index!`$C`$2"
